$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.950.13"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.330.67"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "177.11"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.32%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "3.327.17"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +5.73%  "
$ws.Range("E11").Value = "  +1.69%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "47.08"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.93%  "
$ws.Range("E13").Value = "  +2.20%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "682.96"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.871.15"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("E16").Value = "  +1.89%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "67.956.81"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  -0.41%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.334.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.41%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +3.18%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.897"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("E24").Value = "  +0.31%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "99.45"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.15%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.54"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.10%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "33.13"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("E31").Value = "  +6.08%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "566.56"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  +1.96%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.93%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "57.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("E36").Value = "  -0.19%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.704.48"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("E38").Value = "  +2.54%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "34.56"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +9.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.133"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -2.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "130.16"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
